$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 17, keeping only the header row and row 2
$ws.Range("A3:B17").EntireRow.Delete() | Out-Null

# Update row 2 values
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 11.51866285751828
